$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 134.33333
$ws.Range("I9").Value = 39
$ws.Range("K9").Value = 39
$ws.Range("M9").Value = 130

# Row 70
$ws.Range("H70").Value = 4221.65
$ws.Range("I70").Value = 3395
$ws.Range("J70").Value = 4313.5
$ws.Range("K70").Value = 10185
$ws.Range("L70").Value = 12940.5
$ws.Range("M70").Value = -9915
$ws.Range("N70").Value = -13480.5

# Row 73
$ws.Range("H73").Value = 4221.65
$ws.Range("I73").Value = 3395
$ws.Range("J73").Value = 4313.5
$ws.Range("K73").Value = 10185
$ws.Range("L73").Value = 12940.5
$ws.Range("M73").Value = -9249
$ws.Range("N73").Value = -14812.5

# Row 76
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 5000
$ws.Range("M76").Value = -4685

# Row 79
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 5000
$ws.Range("M79").Value = -3908

# Row 80
$ws.Range("H80").Value = 801.1111
$ws.Range("I80").Value = 801.1111
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2403.3333
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -1405.3333

# Row 83
$ws.Range("H83").Value = 801.1111
$ws.Range("I83").Value = 801.1111
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7209.9999
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -2217.9999

# Row 98
$ws.Range("H98").Value = 1828.6471
$ws.Range("I98").Value = 1639.6364
$ws.Range("K98").Value = 1639.6364
$ws.Range("M98").Value = -141.6364000000001

# Row 100
$ws.Range("H100").Value = 5832.3335
$ws.Range("I100").Value = 4498.75
$ws.Range("K100").Value = 4498.75
$ws.Range("M100").Value = -3957.75

# Row 112
$ws.Range("H112").Value = 1715.2759
$ws.Range("I112").Value = 1399.6
$ws.Range("J112").Value = 1781.0416
$ws.Range("K112").Value = 4198.799999999999
$ws.Range("L112").Value = 5343.1248
$ws.Range("M112").Value = -3090.799999999999
$ws.Range("N112").Value = -7559.1248

# Row 113
$ws.Range("H113").Value = 4319
$ws.Range("I113").Value = 4056.2856
$ws.Range("J113").Value = 4686.8
$ws.Range("K113").Value = 4056.2856
$ws.Range("L113").Value = 4686.8
$ws.Range("M113").Value = -802.2856000000002
$ws.Range("N113").Value = -11194.8

# Row 122
$ws.Range("H122").Value = 1828.6471
$ws.Range("I122").Value = 1639.6364
$ws.Range("K122").Value = 4918.9092
$ws.Range("M122").Value = -2468.9092

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 667.3
$ws.Range("I4").Value = 667.3
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 667.3
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -551.3

# Row 122
$ws.Range("H122").Value = 2957.5
$ws.Range("I122").Value = 2957.5
$ws.Range("K122").Value = 8872.5
$ws.Range("M122").Value = -6422.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3038.5
$ws.Range("I86").Value = 2123
$ws.Range("K86").Value = 2123
$ws.Range("M86").Value = -1000

# Row 89
$ws.Range("H89").Value = 3038.5
$ws.Range("I89").Value = 2123
$ws.Range("K89").Value = 10615
$ws.Range("M89").Value = -4999

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 2000715
$ws.Range("I3").Value = 3333850
$ws.Range("K3").Value = 3333850
$ws.Range("M3").Value = -3333737

# Row 99
$ws.Range("H99").Value = 2736.6667
$ws.Range("I99").Value = 2736.6667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2736.6667
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -1238.6667

# Row 126
$ws.Range("H126").Value = 2736.6667
$ws.Range("I126").Value = 2736.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8210.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5740.000100000001
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 293.1875
$ws.Range("I2").Value = 51.090908
$ws.Range("J2").Value = 825.8
$ws.Range("K2").Value = 306.545448
$ws.Range("L2").Value = 4954.799999999999
$ws.Range("M2").Value = -193.545448
$ws.Range("N2").Value = -5180.799999999999

# Row 14
$ws.Range("H14").Value = 324.66666
$ws.Range("I14").Value = 324.66666
$ws.Range("K14").Value = 973.9999799999999
$ws.Range("M14").Value = -800.9999799999999

# Row 92
$ws.Range("H92").Value = 875.25
$ws.Range("I92").Value = 833.6667
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 2501.0001
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -1253.0001
$ws.Range("N92").Value = -5496

# Row 109
$ws.Range("H109").Value = 1869.5555
$ws.Range("I109").Value = 1831.5
$ws.Range("K109").Value = 5494.5
$ws.Range("M109").Value = -4454.5

# Row 118
$ws.Range("H118").Value = 1500
$ws.Range("I118").Value = 1500
$ws.Range("K118").Value = 4500
$ws.Range("M118").Value = -3257

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 457.95456
$ws.Range("I2").Value = 477.05
$ws.Range("J2").Value = 267
$ws.Range("K2").Value = 477.05
$ws.Range("L2").Value = 267
$ws.Range("M2").Value = -364.05
$ws.Range("N2").Value = -493

# Row 3
$ws.Range("H3").Value = 380
$ws.Range("I3").Value = 350
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 350
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = -234
$ws.Range("N3").Value = -732

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 153.125
$ws.Range("I55").Value = 159.66667
$ws.Range("K55").Value = 159.66667
$ws.Range("M55").Value = 13.33332999999999

# Row 82
$ws.Range("H82").Value = 1138.7858
$ws.Range("I82").Value = 541.25
$ws.Range("J82").Value = 1377.8
$ws.Range("K82").Value = 541.25
$ws.Range("L82").Value = 1377.8
$ws.Range("M82").Value = -180.25
$ws.Range("N82").Value = -2099.8

# Row 85
$ws.Range("H85").Value = 1138.7858
$ws.Range("I85").Value = 541.25
$ws.Range("J85").Value = 1377.8
$ws.Range("K85").Value = 541.25
$ws.Range("L85").Value = 1377.8
$ws.Range("M85").Value = 706.75
$ws.Range("N85").Value = -3873.8

$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 10000
$ws.Range("I10").Value = 10000
$ws.Range("K10").Value = 10000
$ws.Range("M10").Value = -9831

# Row 107
$ws.Range("H107").Value = 411.125
$ws.Range("I107").Value = 399.75
$ws.Range("J107").Value = 422.5
$ws.Range("K107").Value = 1199.25
$ws.Range("L107").Value = 1267.5
$ws.Range("M107").Value = 720.75
$ws.Range("N107").Value = -5107.5
